# Add "Building" column to each of the 3 Effect blocks (Effect 1, Effect 2,
# Effect 3) between the "Amount" and "Job" columns, mirroring the existing
# "Trigger Building" column. This supports achievement gating that checks a
# building id as part of an effect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Effect 1 block: insert new column before current column S ("Effect 1 Job").
$ws.Range("S1").EntireColumn.Insert()

# Effect 2 block: after the first insert, the old "Effect 2 Job" column
# (originally AB) now sits at AC. Insert the new Building column before it.
$ws.Range("AC1").EntireColumn.Insert()

# Effect 3 block: after the first two inserts, the old "Effect 3 Job" column
# (originally AK) now sits at AM. Insert the new Building column before it.
$ws.Range("AM1").EntireColumn.Insert()

# Populate the header row, type row, and key row for the three new columns.
$ws.Range("S1").Value = "Effect 1 Building"
$ws.Range("S2").Value = "uint?"
$ws.Range("S3").Value = "building"

$ws.Range("AC1").Value = "Effect 2 Building"
$ws.Range("AC2").Value = "uint?"
$ws.Range("AC3").Value = "building"

$ws.Range("AM1").Value = "Effect 3 Building"
$ws.Range("AM2").Value = "uint?"
$ws.Range("AM3").Value = "building"
